# Update cryptocurrency Price (D) and Volume(1h) (E) figures to refreshed
# values scraped on Fri Feb  3 18:56:47 UTC 2023.
#
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the source data, which are text cells) instead
# of auto-converting numeric-looking strings into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'332.79"
$ws.Range("E2").Value = "'0.65%"

# Row 3
$ws.Range("D3").Value = "'41.28"
$ws.Range("E3").Value = "'1.30%"

# Row 4
$ws.Range("D4").Value = "'5.689"
$ws.Range("E4").Value = "'-6.64%"

# Row 5
$ws.Range("D5").Value = "'0.08079"
$ws.Range("E5").Value = "'-1.49%"

# Row 6
$ws.Range("D6").Value = "'2.035"
$ws.Range("E6").Value = "'2.98%"

# Row 7
$ws.Range("D7").Value = "'8.734"
$ws.Range("E7").Value = "'-0.94%"

# Row 8
$ws.Range("D8").Value = "'4.534"
$ws.Range("E8").Value = "'-1.04%"

# Row 9
$ws.Range("D9").Value = "'3.000"
$ws.Range("E9").Value = "'1.96%"

# Row 10
$ws.Range("D10").Value = "'0.9219"
$ws.Range("E10").Value = "'-3.12%"

# Row 11
$ws.Range("E11").Value = "'-7.19%"

# Row 12
$ws.Range("D12").Value = "'0.1946"
$ws.Range("E12").Value = "'-2.96%"

# Row 13
$ws.Range("D13").Value = "'8.824"
$ws.Range("E13").Value = "'-16.41%"

# Row 14
$ws.Range("D14").Value = "'0.09364"
$ws.Range("E14").Value = "'0.81%"

# Row 15
$ws.Range("D15").Value = "'0.03699"
$ws.Range("E15").Value = "'5.60%"

# Row 16
$ws.Range("D16").Value = "'0.1052"
$ws.Range("E16").Value = "'9.05%"

# Row 17
$ws.Range("D17").Value = "'0.001297"
$ws.Range("E17").Value = "'-1.08%"

# Row 18
$ws.Range("D18").Value = "'0.006303"
$ws.Range("E18").Value = "'2.77%"

# Row 19
$ws.Range("D19").Value = "'3.368"
$ws.Range("E19").Value = "'0.46%"

# Row 20
$ws.Range("E20").Value = "'-1.61%"

# Row 21
$ws.Range("D21").Value = "'0.1422"
$ws.Range("E21").Value = "'-1.09%"

# Row 22
$ws.Range("D22").Value = "'0.2660"
$ws.Range("E22").Value = "'9.34%"

# Row 23
$ws.Range("D23").Value = "'0.04422"
$ws.Range("E23").Value = "'0.00%"

# Row 24
$ws.Range("D24").Value = "'0.001260"
$ws.Range("E24").Value = "'0.48%"

# Row 25
$ws.Range("D25").Value = "'0.004321"
$ws.Range("E25").Value = "'-2.56%"

# Row 26
$ws.Range("D26").Value = "'0.0001242"
$ws.Range("E26").Value = "'14.28%"

# Row 39
$ws.Range("D39").Value = "'0.02878"
$ws.Range("E39").Value = "'14.19%"

# Row 40
$ws.Range("E40").Value = "'3.40%"

# Row 41
$ws.Range("D41").Value = "'0.007756"
$ws.Range("E41").Value = "'3.53%"

# Row 42
$ws.Range("D42").Value = "'0.009934"
$ws.Range("E42").Value = "'9.97%"

# Row 43
$ws.Range("D43").Value = "'0.1419"

# Row 44
$ws.Range("D44").Value = "'0.002234"
$ws.Range("E44").Value = "'9.19%"

# Row 45
$ws.Range("D45").Value = "'0.01104"
$ws.Range("E45").Value = "'3.95%"

# Row 46
$ws.Range("D46").Value = "'0.00006814"
$ws.Range("E46").Value = "'0.07%"

# Row 47
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.31%"

# Row 48
$ws.Range("D48").Value = "'0.002282"
$ws.Range("E48").Value = "'27.01%"

# Row 49
$ws.Range("D49").Value = "'0.003014"
$ws.Range("E49").Value = "'-13.18%"

# Row 50
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.31%"

# Row 51
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.31%"
